$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "donc1236"
$ws.Range("E3").Value = "donc4566"
$ws.Range("E4").Value = "donc7896"

$ws.Range("E4").Select()
